$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '62.660.77'
Set-TextValue "E2" '  -0.70%  '

Set-TextValue "D3" '2.453.87'
Set-TextValue "E3" '  -0.71%  '

Set-TextValue "E4" '  +0.03%  '

Set-TextValue "D5" '570.66'
Set-TextValue "E5" '  -1.26%  '

Set-TextValue "D6" '146.16'
Set-TextValue "E6" '  -0.48%  '

Set-TextValue "E7" '  -0.02%  '

Set-TextValue "E8" '  -1.97%  '

Set-TextValue "E9" '  -1.41%  '

Set-TextValue "E10" '  -0.35%  '

Set-TextValue "D11" '5.17'
Set-TextValue "E11" '  -2.11%  '

Set-TextValue "E12" '  -1.38%  '

Set-TextValue "D13" '28.58'
Set-TextValue "E13" '  -1.28%  '

Set-TextValue "E14" '  -3.63%  '

Set-TextValue "D15" '2.899.35'
Set-TextValue "E15" '  -0.69%  '

Set-TextValue "D16" '62.490.60'
Set-TextValue "E16" '  -0.81%  '

Set-TextValue "D17" '2.475.19'
Set-TextValue "E17" '  +0.14%  '

Set-TextValue "D18" '7.66'
Set-TextValue "E18" '  -6.56%  '

Set-TextValue "E19" '  -3.23%  '

Set-TextValue "D20" '2.23'
Set-TextValue "E20" '  -0.44%  '

Set-TextValue "D21" '321.27'
Set-TextValue "E21" '  -2.41%  '

Set-TextValue "D22" '4.13'
Set-TextValue "E22" '  -0.18%  '

Set-TextValue "D23" '1.00'
Set-TextValue "E23" '  +0.03%  '

Set-TextValue "D24" '9.91'
Set-TextValue "E24" '  +2.76%  '

Set-TextValue "D25" '64.65'
Set-TextValue "E25" '  -2.46%  '

Set-TextValue "D26" '646.28'
Set-TextValue "E26" '  -3.51%  '

Set-TextValue "E27" '  -0.52%  '

Set-TextValue "B28" 'PEPE'
Set-TextValue "C28" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D28" '0.0₃0951'
Set-TextValue "E28" '  -4.14%  '

Set-TextValue "B29" 'Binance-PegBSC-USD'
Set-TextValue "C29" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D29" '0.996'
Set-TextValue "E29" '  -0.56%  '

Set-TextValue "E30" '  -3.60%  '

Set-TextValue "D31" '7.82'
Set-TextValue "E31" '  -2.85%  '

Set-TextValue "E32" '  -3.36%  '

Set-TextValue "E33" '  -0.43%  '

Set-TextValue "E34" '  -0.07%  '

Set-TextValue "D35" '1.48'
Set-TextValue "E35" '  -4.26%  '

Set-TextValue "D36" '4.63'
Set-TextValue "E36" '  -3.16%  '

Set-TextValue "D37" '150.66'
Set-TextValue "E37" '  -1.15%  '

Set-TextValue "D38" '18.53'
Set-TextValue "E38" '  -1.34%  '

Set-TextValue "D39" '0.363'
Set-TextValue "E39" '  -2.50%  '

Set-TextValue "D40" '5.31'
Set-TextValue "E40" '  -2.81%  '

Set-TextValue "E41" '  -4.15%  '

Set-TextValue "E42" '  -3.80%  '

Set-TextValue "B43" 'BabyDogeCoin'
Set-TextValue "C43" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D43" '0.0₆0312'
Set-TextValue "E43" '  +2.18%  '

Set-TextValue "B44" 'USDe'
Set-TextValue "C44" 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue "D44" '1.01'
Set-TextValue "E44" '  +0.62%  '

Set-TextValue "D45" '152.47'
Set-TextValue "E45" '  +0.43%  '

Set-TextValue "E46" '  +1.76%  '

Set-TextValue "E47" '  -2.14%  '

Set-TextValue "D48" '0.602'
Set-TextValue "E48" '  -0.61%  '

Set-TextValue "D49" '19.92'
Set-TextValue "E49" '  -3.45%  '

Set-TextValue "D50" '0.0504'
Set-TextValue "E50" '  -1.59%  '

Set-TextValue "E51" '  -2.03%  '
